# edit.ps1
# Applies the "Updated symbol list" price/volume refresh to the crypto tracker sheet.
# For each affected cell we explicitly mark the cell as Text ("@" number format)
# before writing the value, so that Excel does not auto-convert numeric-looking
# strings like "310.44" or percentages like "1.76%" into actual number/percentage
# values. This preserves the original text-cell semantics used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.76%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.15%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.105"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.41%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08198"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.55%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.056"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.38%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.948"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.958"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "11.15%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9256"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.10%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1121"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "14.95%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1916"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.48%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09297"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.12%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03650"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.58%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09907"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.20%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001431"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.78%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005849"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.40%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.471"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.20%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.123"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.37%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3395"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.81%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1308"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.96%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.095"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.39%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2211"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.40%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04529"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.94%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001226"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.54%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004809"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.21%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001250"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.78%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004442"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.22%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01968"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.00%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04892"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.14%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007610"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.22%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009053"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "17.20%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1384"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.91%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002188"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.15%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01163"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.50%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006550"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.13%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.08%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "179.94"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "247.05%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001498"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-20.95%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.08%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.08%"
